$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "58.452.09", "  -1.62%  "),
    @("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "2.620.32", "  +0.65%  "),
    @("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  -0.02%  "),
    @("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "534.58", "  -0.53%  "),
    @("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "141.97", "  +0.38%  "),
    @("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "0.999", "  -0.03%  "),
    @("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.566", "  +0.08%  "),
    @("LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "2.627.60", "  +0.78%  "),
    @("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "6.97", "  +8.23%  "),
    @("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.100", "  -1.98%  "),
    @("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.335", "  -0.27%  "),
    @("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.134", "  +0.40%  "),
    @("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "3.087.20", "  +0.70%  "),
    @("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "58.407.90", "  -1.60%  "),
    @("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "20.82", "  +1.04%  "),
    @("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "2.607.30", "  +0.36%  "),
    @("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000132", "  -0.79%  "),
    @("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "334.65", "  -1.89%  "),
    @("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "4.37", "  +0.22%  "),
    @("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "10.13", "  +0.22%  "),
    @("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "6.20", "  -2.71%  "),
    @("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.999", "  -0.11%  "),
    @("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "66.57", "  -1.40%  "),
    @("Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.415", "  +1.77%  "),
    @("Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.163", "  -0.70%  "),
    @("Binance-PegBSC-USD", "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd", "0.999", "  -0.11%  "),
    @("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "7.10", "  -1.64%  "),
    @("USDe", "https://coinranking.com/coin/exbfr2U-0+usde-usde", "0.998", "  -0.07%  "),
    @("PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0₃0731", "  -2.02%  "),
    @("PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.64", "  -2.21%  "),
    @("Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "5.80", "  +0.01%  "),
    @("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "150.64", "  +0.09%  "),
    @("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "18.70", "  -0.70%  "),
    @("NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "3.89", "  -2.01%  "),
    @("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "1.10", "  -1.40%  "),
    @("SuiNetwork", "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui", "0.819", "  -1.78%  "),
    @("Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "1.42", "  -3.15%  "),
    @("Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "0.809", "  -2.31%  "),
    @("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "3.58", "  +0.93%  "),
    @("Bittensor", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao", "281.77", "  +2.35%  "),
    @("FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "0.999", "  -0.04%  "),
    @("Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.593", "  -0.53%  "),
    @("WhiteBITCoin", "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt", "10.67", "  -0.53%  "),
    @("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.0529", "  +1.30%  "),
    @("Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.0935", "  -2.10%  "),
    @("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "18.89", "  +1.93%  "),
    @("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0225", "  +0.96%  "),
    @("Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "1.936.33", "  -0.74%  "),
    @("RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "4.43", "  -1.78%  "),
    @("InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "17.83", "  -3.45%  "),
)

$numericPattern = '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$'

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]

    $dCell = $ws.Cells.Item($r, 4)
    $dVal = $row[2]
    if ($dVal -match $numericPattern) {
        $dCell.NumberFormat = "@"
    }
    $dCell.Value = $dVal

    $ws.Cells.Item($r, 5).Value = $row[3]
}
